$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so numeric-looking values (e.g. "1.001") are not
# auto-converted to numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.273.05"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "1.893.23"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "323.02"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.5177"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "0.08417"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "42.70"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "23.06"
$ws.Range("E12").Value = "  +10.25%  "
$ws.Range("D13").Value = "6.432"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "1.887.69"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "7.312"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "94.23"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "0.00001108"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "0.06649"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "30.261.31"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "2.228"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "2.106.00"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "21.54"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").Value = "162.12"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "2.332"
$ws.Range("E29").Value = "  -5.47%  "
$ws.Range("D30").Value = "129.24"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "6.101"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "3.745"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "0.02493"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").Value = "0.06539"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "5.333"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "0.2197"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "8.808"
$ws.Range("E40").Value = "  -3.91%  "
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "0.6502"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "1.227"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "0.6085"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "13.27"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "3.681"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "2.053"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "1.235"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "124.46"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").Value = "79.07"
$ws.Range("E51").Value = "  +0.78%  "
